$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NPCs")

# Add the new NPC row (row 34) - "Federation Chancellor" quest NPC
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = "FederationChancellor"
$ws.Cells.Item(34, 3).Value = "Federation Chancellor"
$ws.Cells.Item(34, 4).Value = 2
$ws.Cells.Item(34, 5).Value = "Delusional Memories"
$ws.Cells.Item(34, 9).Value = 80
$ws.Cells.Item(34, 10).Value = 96
